$d = $word.ActiveDocument

# --- Locate the paragraph that holds "Dia 16/10: 2hr 23min (1 dia)" ---
$oldText = "Dia 16/10: 2hr 23min (1 dia)"
$findRange = $d.Content
$found = $findRange.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                                  $true, 1, $false, $null, 0)
if (-not $found) {
    throw "Could not find target paragraph text: $oldText"
}

# Find the Paragraph object whose Range starts at the same position as the match.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Start -eq $findRange.Start) {
        $target = $cand
        break
    }
}
if ($null -eq $target) {
    throw "Could not resolve paragraph object for match"
}

$paraRange = $target.Range

# --- Capture the paragraph's own formatting so we can rebuild it faithfully ---
$fontName = $paraRange.Font.Name
$fontSize = $paraRange.Font.Size

# Pull the exact <w:p ...> opening tag (with its w14:paraId / rsid attributes, etc.)
# straight out of the live package so none of that metadata is lost.
$openXml = $paraRange.WordOpenXML
$pTagOpen = "<w:p>"
if ($openXml -match '(<w:p [^>]*>)') {
    $pTagOpen = $matches[1]
}

# --- Build the replacement run-property block matching the paragraph's font ---
$rPr = "<w:rPr><w:rFonts w:ascii=""$fontName"" w:hAnsi=""$fontName"" w:cs=""$fontName""/><w:sz w:val=""$([int]($fontSize*2))""/><w:szCs w:val=""$([int]($fontSize*2))""/></w:rPr>"

# --- The five runs that replace the single original run ---
$runs =
    "<w:r>$rPr<w:t xml:space=""preserve"">Dia 16/10: </w:t></w:r>" +
    "<w:r>$rPr<w:t>3</w:t></w:r>" +
    "<w:r>$rPr<w:t xml:space=""preserve"">hr </w:t></w:r>" +
    "<w:r>$rPr<w:t>37</w:t></w:r>" +
    "<w:r>$rPr<w:t>min (1 dia)</w:t></w:r>"

$pPr = "<w:pPr><w:spacing w:line=""360"" w:lineRule=""auto""/><w:jc w:val=""both""/>$rPr</w:pPr>"

$newParaXml = "$pTagOpen$pPr$runs</w:p>"
$wrapped = "<pkg:xml xmlns:w=""http://schemas.openxmlformats.org/wordprocessingml/2006/main"" " +
           "xmlns:w14=""http://schemas.microsoft.com/office/word/2010/wordml"">$newParaXml</pkg:xml>"

# Replace the whole paragraph (including its end-of-paragraph mark) in one shot.
$paraRange.InsertXML($newParaXml)
